$d = $word.ActiveDocument

# ---------------------------------------------------------------------
# Edit 1: "The project manager and team leader will review the change
# request." -> split into three runs, replacing "team leader " with
# "product owner ".
# ---------------------------------------------------------------------
$search1 = "The project manager and team leader will review the change request."
$r1 = $d.Content
$found1 = $r1.Find.Execute($search1, $false, $false, $false, $false, $false, $true, 1, $false, "", 0)
if (-not $found1) {
    throw "Could not find target sentence for edit 1"
}

$base1 = $r1.Start

$part1a = "The project manager and "
$part1b = "product owner "
$part1c = "will review the change request."

$r1.Text = $part1a + $part1b + $part1c

$p1aEnd = $base1 + $part1a.Length
$p1bEnd = $p1aEnd + $part1b.Length

# Toggling a character formatting property on/off over the middle
# span's actual (non-empty) range forces the run to be split off from
# its neighbours, while leaving the effective formatting unchanged
# (same rPr as before/after, matching the target OOXML).
$rMid1 = $d.Range($p1aEnd, $p1bEnd)
$rMid1.Bold = 1
$rMid1.Bold = 0

# ---------------------------------------------------------------------
# Edit 2: "A change request log will be submitted by a team member or
# project sponsor and will be conducted for impact analysis. " ->
# split into three runs so that "ject sponsor " becomes "duct owner ".
# ---------------------------------------------------------------------
$search2 = "A change request log will be submitted by a team member or project sponsor and will be conducted for impact analysis. "
$r2 = $d.Content
$found2 = $r2.Find.Execute($search2, $false, $false, $false, $false, $false, $true, 1, $false, "", 0)
if (-not $found2) {
    throw "Could not find target sentence for edit 2"
}

$base2 = $r2.Start

$part2a = "A change request log will be submitted by a team member or pro"
$part2b = "duct owner "
$part2c = "and will be conducted for impact analysis. "

$r2.Text = $part2a + $part2b + $part2c

$p2aEnd = $base2 + $part2a.Length
$p2bEnd = $p2aEnd + $part2b.Length

$rMid2 = $d.Range($p2aEnd, $p2bEnd)
$rMid2.Bold = 1
$rMid2.Bold = 0
